$d = $word.ActiveDocument

# --- Paragraph 2: "Making a bet of 5 with a a balance of 5:" ---
# Collapse the three runs (with the gramStart/gramEnd proofErr markers around
# the stray "a") into a single run, text unchanged.
$r2 = $d.Content
[void]$r2.Find.Execute("Making a bet of 5 with a a balance of 5:", $true, $false, $false, $false, $false, $true, 1, $false, "Making a bet of 5 with a a balance of 5:", 2)

# --- Paragraph 1: "Bug:" -> "Bug" / bookmark / "2" / ":" ---
$r1 = $d.Content
[void]$r1.Find.Execute("Bug:")
$bugStart = $r1.Start
$bugEnd = $r1.End

# Type "2" right before the colon -> "Bug2:"
$insPos = $d.Range($bugEnd - 1, $bugEnd - 1)
$insPos.InsertBefore("2")

# Break "Bug" into its own run (still carrying the preceding <w:br/>) by
# toggling a format and reverting it - this forces a run boundary without
# altering the visible formatting.
$bugRange = $d.Range($bugStart - 1, $bugEnd - 1)
$bugRange.Bold = 1
$bugRange.Bold = 0

# Break "2" away from ":" the same way.
$twoRange = $d.Range($bugEnd - 1, $bugEnd)
$twoRange.Bold = 1
$twoRange.Bold = 0

# Drop the _GoBack bookmark at the edit point (between "Bug" and "2") - this
# also removes it from its old location since it's a singleton bookmark.
$bmPos = $d.Range($bugEnd - 1, $bugEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmPos)
